# Update "想去人数" (want-to-go count) values by +1 (or +2 for row 14/18 on 展览/全部类型)
# as published by the gh-pages data refresh at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(3, 6).Value = 782
$ws1.Cells.Item(4, 6).Value = 625
$ws1.Cells.Item(12, 6).Value = 712
$ws1.Cells.Item(13, 6).Value = 1205
$ws1.Cells.Item(14, 6).Value = 238
$ws1.Cells.Item(15, 6).Value = 970
$ws1.Cells.Item(18, 6).Value = 19
$ws1.Cells.Item(19, 6).Value = 349
$ws1.Cells.Item(22, 6).Value = 498
$ws1.Cells.Item(23, 6).Value = 521
$ws1.Cells.Item(24, 6).Value = 710

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(10, 6).Value = 461

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2, 6).Value = 136

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(3, 6).Value = 782
$ws4.Cells.Item(5, 6).Value = 625
$ws4.Cells.Item(7, 6).Value = 136
$ws4.Cells.Item(16, 6).Value = 712
$ws4.Cells.Item(17, 6).Value = 1205
$ws4.Cells.Item(18, 6).Value = 238
$ws4.Cells.Item(19, 6).Value = 970
$ws4.Cells.Item(22, 6).Value = 19
$ws4.Cells.Item(23, 6).Value = 349
$ws4.Cells.Item(31, 6).Value = 498
$ws4.Cells.Item(32, 6).Value = 521
$ws4.Cells.Item(33, 6).Value = 710
$ws4.Cells.Item(37, 6).Value = 461
